{"js": "// Update the date stamp and every \"A\u00d7B=C\" answer cell to the new values\n// from the commit (output regenerated at c986bee). Each old value is\n// unique in the document, so a plain text search-and-replace for every\n// pair is sufficient and keeps all run/paragraph formatting untouched.\nconst replacements = [\n  [\"2024-10-02 Wednesday\", \"2024-10-03 Thursday\"],\n  [\"627\u00d73=1881\", \"773\u00d77=5411\"],\n  [\"125\u00d73=375\", \"776\u00d72=1552\"],\n  [\"526\u00d78=4208\", \"916\u00d74=3664\"],\n  [\"291\u00d72=582\", \"738\u00d79=6642\"],\n  [\"335\u00d74=1340\", \"762\u00d74=3048\"],\n  [\"304\u00d77=2128\", \"964\u00d72=1928\"],\n  [\"363\u00d73=1089\", \"432\u00d78=3456\"],\n  [\"607\u00d72=1214\", \"184\u00d79=1656\"],\n  [\"971\u00d75=4855\", \"504\u00d78=4032\"],\n  [\"687\u00d73=2061\", \"624\u00d75=3120\"],\n  [\"364\u00d73=1092\", \"695\u00d74=2780\"],\n  [\"448\u00d78=3584\", \"370\u00d77=2590\"],\n  [\"379\u00d78=3032\", \"314\u00d79=2826\"],\n  [\"499\u00d77=3493\", \"780\u00d78=6240\"],\n  [\"744\u00d77=5208\", \"406\u00d74=1624\"],\n  [\"963\u00d77=6741\", \"367\u00d77=2569\"],\n  [\"852\u00d72=1704\", \"278\u00d78=2224\"],\n  [\"154\u00d75=770\", \"147\u00d76=882\"],\n  [\"246\u00d73=738\", \"619\u00d79=5571\"],\n  [\"949\u00d72=1898\", \"563\u00d73=1689\"],\n  [\"408\u00d77=2856\", \"362\u00d79=3258\"],\n  [\"390\u00d73=1170\", \"814\u00d73=2442\"],\n  [\"669\u00d72=1338\", \"603\u00d72=1206\"],\n  [\"732\u00d72=1464\", \"396\u00d72=792\"],\n  [\"854\u00d79=7686\", \"631\u00d73=1893\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date stamp and every \"A\u00d7B=C\" answer cell to the new values\n# from the commit (output regenerated at c986bee). Each old value is\n# unique in the document, so a plain Find/Replace for every pair is\n# sufficient and leaves all other run/paragraph formatting untouched.\n\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = [ordered]@{\n    \"2024-10-02 Wednesday\" = \"2024-10-03 Thursday\"\n    \"627\u00d73=1881\"           = \"773\u00d77=5411\"\n    \"125\u00d73=375\"            = \"776\u00d72=1552\"\n    \"526\u00d78=4208\"           = \"916\u00d74=3664\"\n    \"291\u00d72=582\"            = \"738\u00d79=6642\"\n    \"335\u00d74=1340\"           = \"762\u00d74=3048\"\n    \"304\u00d77=2128\"           = \"964\u00d72=1928\"\n    \"363\u00d73=1089\"           = \"432\u00d78=3456\"\n    \"607\u00d72=1214\"           = \"184\u00d79=1656\"\n    \"971\u00d75=4855\"           = \"504\u00d78=4032\"\n    \"687\u00d73=2061\"           = \"624\u00d75=3120\"\n    \"364\u00d73=1092\"           = \"695\u00d74=2780\"\n    \"448\u00d78=3584\"           = \"370\u00d77=2590\"\n    \"379\u00d78=3032\"           = \"314\u00d79=2826\"\n    \"499\u00d77=3493\"           = \"780\u00d78=6240\"\n    \"744\u00d77=5208\"           = \"406\u00d74=1624\"\n    \"963\u00d77=6741\"           = \"367\u00d77=2569\"\n    \"852\u00d72=1704\"           = \"278\u00d78=2224\"\n    \"154\u00d75=770\"            = \"147\u00d76=882\"\n    \"246\u00d73=738\"            = \"619\u00d79=5571\"\n    \"949\u00d72=1898\"           = \"563\u00d73=1689\"\n    \"408\u00d77=2856\"           = \"362\u00d79=3258\"\n    \"390\u00d73=1170\"           = \"814\u00d73=2442\"\n    \"669\u00d72=1338\"           = \"603\u00d72=1206\"\n    \"732\u00d72=1464\"           = \"396\u00d72=792\"\n    \"854\u00d79=7686\"           = \"631\u00d73=1893\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll)\n}\n"}
